$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '60.978.56'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '2.370.23'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '547.76'
$ws.Range('E5').Value = '  -0.36%  '
Set-TextValue $ws.Range('D6') '137.58'
$ws.Range('E6').Value = '  -3.20%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -2.55%  '
$ws.Range('D9').Value = '2.371.51'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('E11').Value = '  +1.12%  '
Set-TextValue $ws.Range('D12') '5.33'
$ws.Range('E12').Value = '  +0.62%  '
Set-TextValue $ws.Range('D13') '0.346'
$ws.Range('E13').Value = '  -0.67%  '
Set-TextValue $ws.Range('D14') '24.94'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').Value = '2.781.01'
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '60.914.26'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '2.388.51'
$ws.Range('E18').Value = '  -0.55%  '
Set-TextValue $ws.Range('D19') '10.77'
$ws.Range('E19').Value = '  -0.59%  '
Set-TextValue $ws.Range('D20') '4.13'
$ws.Range('E20').Value = '  -0.52%  '
Set-TextValue $ws.Range('D21') '319.59'
$ws.Range('E21').Value = '  -0.03%  '
Set-TextValue $ws.Range('D22') '6.67'
$ws.Range('E22').Value = '  -1.03%  '
Set-TextValue $ws.Range('D24') '64.15'
$ws.Range('E24').Value = '  +0.55%  '
Set-TextValue $ws.Range('D25') '1.67'
$ws.Range('E25').Value = '  -13.20%  '
Set-TextValue $ws.Range('D26') '8.46'
$ws.Range('E26').Value = '  +2.48%  '
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '2.485.38'
$ws.Range('E28').Value = '  -1.22%  '
Set-TextValue $ws.Range('D29') '8.11'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D30') '0.150'
$ws.Range('E30').Value = '  +2.15%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D31') '503.71'
$ws.Range('E31').Value = '  -6.03%  '
$ws.Range('D32').Value = '0.0₃0875'
$ws.Range('E32').Value = '  -7.09%  '
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('E35').Value = '  -5.15%  '
Set-TextValue $ws.Range('D36') '0.999'
$ws.Range('E36').Value = '  -0.05%  '
Set-TextValue $ws.Range('D37') '4.65'
$ws.Range('E37').Value = '  -1.93%  '
Set-TextValue $ws.Range('D38') '0.377'
Set-TextValue $ws.Range('D39') '1.86'
$ws.Range('E39').Value = '  +0.34%  '
Set-TextValue $ws.Range('D40') '18.52'
$ws.Range('E40').Value = '  +1.82%  '
Set-TextValue $ws.Range('D41') '5.32'
$ws.Range('E41').Value = '  -4.39%  '
Set-TextValue $ws.Range('D42') '145.75'
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('E43').Value = '  -0.09%  '
Set-TextValue $ws.Range('D44') '41.57'
$ws.Range('E44').Value = '  +3.08%  '
Set-TextValue $ws.Range('D45') '147.16'
$ws.Range('E45').Value = '  +3.98%  '
Set-TextValue $ws.Range('D46') '3.58'
$ws.Range('E46').Value = '  -1.45%  '
Set-TextValue $ws.Range('D47') '2.04'
$ws.Range('E47').Value = '  -6.85%  '
Set-TextValue $ws.Range('D48') '0.0519'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D49') '0.574'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D50') '19.18'
$ws.Range('E50').Value = '  -5.71%  '
Set-TextValue $ws.Range('D51') '0.0909'
$ws.Range('E51').Value = '  -0.19%  '
